# Updated cryptos list on Sat Aug 17 03:45:04 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '58.864.95'
$ws.Range("E2").Value = '  +1.86%  '
$ws.Range("D3").Value = '2.581.53'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.48'
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.49'
$ws.Range("E6").Value = '  -2.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  -0.39%  '
$ws.Range("D9").Value = '2.592.21'
$ws.Range("E10").Value = '  -0.44%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.328'
$ws.Range("E12").Value = '  +1.68%  '
$ws.Range("E13").Value = '  +3.26%  '
$ws.Range("D14").Value = '3.038.12'
$ws.Range("E14").Value = '  +0.65%  '
$ws.Range("D15").Value = '58.951.51'
$ws.Range("E15").Value = '  +1.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.35'
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").Value = '2.594.51'
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000132'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '336.53'
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.50'
$ws.Range("E22").Value = '  +3.21%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.07'
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.00'
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = '0.0₃0721'
$ws.Range("E30").Value = '  -2.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.93'
$ws.Range("E31").Value = '  -5.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.57'
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.62'
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.00'
$ws.Range("E34").Value = '  -0.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.96'
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("E36").Value = '  -2.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.77'
$ws.Range("E37").Value = '  +2.04%  '
$ws.Range("E38").Value = '  +1.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.821'
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.804'
$ws.Range("E40").Value = '  -7.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.50'
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '269.77'
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.587'
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("E47").Value = '  -0.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.33'
$ws.Range("E48").Value = '  -1.91%  '
$ws.Range("D49").Value = '1.961.77'
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0219'
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.49'
$ws.Range("E51").Value = '  -2.02%  '
